# Update "想去人数" (wanted-to-go count, column F) values on the
# "展览" and "全部类型" sheets to reflect the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1874
$ws1.Range("F4").Value  = 168
$ws1.Range("F6").Value  = 2585
$ws1.Range("F8").Value  = 91
$ws1.Range("F11").Value = 530
$ws1.Range("F13").Value = 333
$ws1.Range("F17").Value = 210
$ws1.Range("F23").Value = 1651
$ws1.Range("F24").Value = 27
$ws1.Range("F25").Value = 402
$ws1.Range("F26").Value = 569

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1874
$ws4.Range("F5").Value  = 168
$ws4.Range("F7").Value  = 2585
$ws4.Range("F9").Value  = 91
$ws4.Range("F12").Value = 530
$ws4.Range("F14").Value = 333
$ws4.Range("F18").Value = 210
$ws4.Range("F24").Value = 1651
$ws4.Range("F25").Value = 27
$ws4.Range("F26").Value = 402
$ws4.Range("F27").Value = 569
